$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the values for rows 2-6 (A:B) with the new data
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 135

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 107

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 81

$ws.Range("A5").Value = 0
$ws.Range("B5").Value = 75

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 71

# Remove rows 7-11 which are no longer part of the data
$ws.Range("A7:B11").EntireRow.Delete()
